# Update column G ("K") values on Sheet1 for rows 2-24.
# This mirrors a regeneration of save_data where the K column values
# (previously computed differently, e.g. from a "Strike#" style source)
# are recalculated and rewritten with new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 5
    3  = 3
    4  = 3
    5  = 6
    6  = 5
    7  = 11
    8  = 8
    9  = 8
    10 = 6
    11 = 7
    12 = 5
    13 = 8
    14 = 8
    15 = 6
    16 = 6
    17 = 3
    18 = 4
    19 = 5
    20 = 7
    21 = 5
    22 = 4
    23 = 5
    24 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
